# Weekly price-report update: a new weekly record is inserted as row 144
# ("Ají" / "Americana (o)" / Región Metropolitana, week of 2022-02-11),
# pushing every existing record from the old row 144 down through the old
# row 242 down by one row (old row 242 becomes the new row 243). The
# used range grows from A1:R242 to A1:R243.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 144; this shifts rows 144-242 down to
# 145-243 and carries the D-column (date) number format onto the new row,
# exactly like Excel's own "Insert Rows" command.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly record.
$ws.Range("A144").Value = 9
$ws.Range("B144").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C144").Value = "Metropolitana"
$ws.Range("D144").Value = 44603
$ws.Range("E144").Value = 13
$ws.Range("F144").Value = 100112021
$ws.Range("G144").Value = "Ají"
$ws.Range("H144").Value = "Americana (o)"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 25
$ws.Range("K144").Value = 18000
$ws.Range("L144").Value = 20000
$ws.Range("M144").Value = 18960
$ws.Range("N144").Value = "`$/caja 25 kilos"
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("P144").Value = 758
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
